$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the first data row of the "posts" table (row 34) ---
# post_id
$ws.Range("B34").Value = 1
# title
$ws.Range("C34").Value = "Introduction Class"
# desc
$ws.Range("D34").Value = "Sample description"
# resources_link / class_link (plain text first, hyperlink applied below)
$ws.Range("E34").Value = "www.google.com"
$ws.Range("F34").Value = "www.google.com"
# course_id
$ws.Range("G34").Value = 10001

# Turn E34 / F34 into hyperlinks pointing at www.google.com
$ws.Hyperlinks.Add($ws.Range("E34"), "www.google.com")
$ws.Hyperlinks.Add($ws.Range("F34"), "www.google.com")

# Hyperlinks.Add() stamps its own ad-hoc xf (distinct from the sheet's
# existing "Hyperlink" look used by F5/F6); reapply that exact formatting so
# E34/F34 match the rest of the hyperlink cells on the sheet.
$ws.Range("F5").Copy()
$ws.Range("E34").PasteSpecial(-4122)
$ws.Range("F5").Copy()
$ws.Range("F34").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- View state: scroll position + active selection ---
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("P17").Select()
